$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new result row (row 59) produced by the latest script run.
# Leading apostrophe forces the date-like string to stay text instead of
# being auto-converted to a date serial by Excel; ClearFormats() then
# drops the resulting quote-prefix style so no stray formatting is left
# on the new cell.
$ws.Range("A59").Value = "'2025-04-25"
$ws.Range("A59").ClearFormats()
$ws.Range("B59").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C59").Value = "NA"
$ws.Range("D59").Value = 1
